# Added filtering options for the Component Analysis
# Shifts the forecast-error rows down by one (row N -> row N+1, for rows 2-10)
# and writes a brand-new set of summary statistics into the now-vacated row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for columns B:G, rows 2-10, before overwriting anything.
$values = @{}
for ($r = 2; $r -le 10; $r++) {
    $srcRange = "B" + $r + ":G" + $r
    $values[$r] = $ws.Range($srcRange).Value2
}

# Shift rows 2-10 down into rows 3-11.
for ($r = 10; $r -ge 2; $r--) {
    $destRange = "B" + ($r + 1) + ":G" + ($r + 1)
    $ws.Range($destRange).Value2 = $values[$r]
}

# New data for row 2.
$ws.Range("B2").Value2 = 0.03483647684766893
$ws.Range("C2").Value2 = 0.9165303275553447
$ws.Range("D2").Value2 = 4.351217448857517
$ws.Range("E2").Value2 = 2.085957202067558
$ws.Range("F2").Value2 = 2.108713024898596
$ws.Range("G2").Value2 = 46
